# Auto-generated Excel COM-interop edit script
# Applies the "Atualizacao de bases das ligas" update for Slovenia Prva Liga sheet:
#  - Fixes the NK Domzale / NK Maribor shared-string mix-up (and every cell that
#    referenced those two team names so the displayed team names stay correct)
#  - Re-pairs the odds/result data that had been attached to the wrong match for
#    rows 9 and 10 (match ids 6814328 / 6814330)
#  - Fills in results + closing odds for two already-listed upcoming fixtures
#    (rows 116 and 117)
#  - Appends 5 new fixture rows (118-122)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the 5 new rows (118-122) and give column A / E the same formatting
#        (bold+bordered id column, date-time number format) used by the rest of
#        the table, by copying it down from the last existing data row (117).
$ws.Range("A117").Copy() | Out-Null
$ws.Range("A118:A122").PasteSpecial(-4122) | Out-Null
$ws.Range("E117").Copy() | Out-Null
$ws.Range("E118:E122").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 2. Table of every cell that changes value (row, column index, type, value)
#     type 's' = string (team name / div name / result letter)
#     type 'n' = numeric
$changes = @(
    @(2,7,'s','NK Domzale'),
    @(5,7,'s','NK Maribor'),
    @(9,2,'n',6814330),
    @(9,7,'s','NK Aluminij'),
    @(9,9,'n',0),
    @(9,10,'s','H'),
    @(9,11,'n',1.363),
    @(9,12,'n',4.5),
    @(9,13,'n',7),
    @(9,14,'n',1.4),
    @(9,15,'n',4.5),
    @(9,16,'n',7),
    @(9,17,'n',-1.25),
    @(9,18,'n',1.85),
    @(9,19,'n',1.95),
    @(9,20,'n',2.75),
    @(9,21,'n',1.8),
    @(9,22,'n',2),
    @(9,23,'n',0.3999999999999999),
    @(9,24,'n',-1),
    @(9,27,'n',0.475),
    @(9,28,'n',-1),
    @(9,29,'n',1),
    @(10,2,'n',6814328),
    @(10,7,'s','NK Bravo'),
    @(10,9,'n',1),
    @(10,10,'s','D'),
    @(10,11,'n',2.35),
    @(10,12,'n',3.1),
    @(10,13,'n',2.9),
    @(10,14,'n',2.15),
    @(10,15,'n',3.1),
    @(10,16,'n',3.3),
    @(10,17,'n',-0.25),
    @(10,18,'n',1.925),
    @(10,19,'n',1.875),
    @(10,20,'n',2.25),
    @(10,21,'n',1.95),
    @(10,22,'n',1.85),
    @(10,23,'n',-1),
    @(10,24,'n',2.1),
    @(10,27,'n',0.4375),
    @(10,28,'n',-0.5),
    @(10,29,'n',0.425),
    @(12,7,'s','NK Domzale'),
    @(14,6,'s','NK Maribor'),
    @(20,7,'s','NK Maribor'),
    @(21,6,'s','NK Domzale'),
    @(23,6,'s','NK Maribor'),
    @(23,7,'s','NK Domzale'),
    @(27,6,'s','NK Domzale'),
    @(29,7,'s','NK Maribor'),
    @(33,7,'s','NK Maribor'),
    @(34,7,'s','NK Domzale'),
    @(36,6,'s','NK Domzale'),
    @(40,6,'s','NK Maribor'),
    @(43,6,'s','NK Domzale'),
    @(44,6,'s','NK Maribor'),
    @(46,7,'s','NK Domzale'),
    @(50,7,'s','NK Maribor'),
    @(52,6,'s','NK Maribor'),
    @(53,6,'s','NK Domzale'),
    @(57,6,'s','NK Domzale'),
    @(58,7,'s','NK Maribor'),
    @(61,7,'s','NK Domzale'),
    @(63,6,'s','NK Maribor'),
    @(68,6,'s','NK Domzale'),
    @(68,7,'s','NK Maribor'),
    @(73,7,'s','NK Domzale'),
    @(74,6,'s','NK Maribor'),
    @(76,7,'s','NK Domzale'),
    @(79,7,'s','NK Maribor'),
    @(81,6,'s','NK Domzale'),
    @(83,6,'s','NK Maribor'),
    @(85,7,'s','NK Domzale'),
    @(86,7,'s','NK Maribor'),
    @(92,7,'s','NK Domzale'),
    @(93,7,'s','NK Maribor'),
    @(96,6,'s','NK Domzale'),
    @(99,6,'s','NK Maribor'),
    @(101,7,'s','NK Domzale'),
    @(105,7,'s','NK Domzale'),
    @(107,6,'s','NK Maribor'),
    @(109,7,'s','NK Maribor'),
    @(110,6,'s','NK Domzale'),
    @(115,6,'s','NK Maribor'),
    @(115,7,'s','NK Domzale'),
    @(116,8,'n',1),
    @(116,9,'n',2),
    @(116,10,'s','A'),
    @(116,14,'n',1.65),
    @(116,15,'n',3.75),
    @(116,18,'n',1.825),
    @(116,19,'n',1.975),
    @(116,20,'n',2.5),
    @(116,21,'n',1.8),
    @(116,22,'n',2),
    @(116,23,'n',-1),
    @(116,24,'n',-1),
    @(116,25,'n',3.75),
    @(116,26,'n',-1),
    @(116,27,'n',0.9750000000000001),
    @(116,28,'n',0.8),
    @(116,29,'n',-1),
    @(117,8,'n',0),
    @(117,9,'n',1),
    @(117,10,'s','A'),
    @(117,21,'n',2.025),
    @(117,22,'n',1.775),
    @(117,23,'n',-1),
    @(117,24,'n',-1),
    @(117,25,'n',1.55),
    @(117,26,'n',-1),
    @(117,27,'n',0.7749999999999999),
    @(117,28,'n',-1),
    @(117,29,'n',0.7749999999999999),
    @(118,1,'n',116),
    @(118,2,'n',6814421),
    @(118,3,'s','Slovenia Prva Liga'),
    @(118,4,'s','Slovenia Prva Liga'),
    @(118,5,'n',45360.45833333334),
    @(118,6,'s','NK Bravo'),
    @(118,7,'s','NK Aluminij'),
    @(118,11,'n',1.666),
    @(118,12,'n',3.5),
    @(118,13,'n',5),
    @(118,14,'n',1.65),
    @(118,15,'n',3.5),
    @(118,16,'n',5.25),
    @(118,17,'n',-0.75),
    @(118,18,'n',1.85),
    @(118,19,'n',1.95),
    @(118,20,'n',2.5),
    @(118,21,'n',2),
    @(118,22,'n',1.8),
    @(118,23,'n',0),
    @(118,24,'n',0),
    @(118,25,'n',0),
    @(118,26,'n',0),
    @(118,27,'n',0),
    @(119,1,'n',117),
    @(119,2,'n',6814423),
    @(119,3,'s','Slovenia Prva Liga'),
    @(119,4,'s','Slovenia Prva Liga'),
    @(119,5,'n',45360.5625),
    @(119,6,'s','NK Celje'),
    @(119,7,'s','NK Maribor'),
    @(119,11,'n',2),
    @(119,12,'n',3.3),
    @(119,13,'n',3.5),
    @(119,14,'n',1.95),
    @(119,15,'n',3.3),
    @(119,16,'n',3.6),
    @(119,17,'n',-0.5),
    @(119,18,'n',2),
    @(119,19,'n',1.8),
    @(119,20,'n',2.5),
    @(119,21,'n',1.925),
    @(119,22,'n',1.875),
    @(119,23,'n',0),
    @(119,24,'n',0),
    @(119,25,'n',0),
    @(119,26,'n',0),
    @(119,27,'n',0),
    @(120,1,'n',118),
    @(120,2,'n',6816449),
    @(120,3,'s','Slovenia Prva Liga'),
    @(120,4,'s','Slovenia Prva Liga'),
    @(120,5,'n',45361.375),
    @(120,6,'s','NK Rogaska'),
    @(120,7,'s','NK Radomlje'),
    @(120,11,'n',2.625),
    @(120,12,'n',3.2),
    @(120,13,'n',2.5),
    @(120,14,'n',2.7),
    @(120,15,'n',3.2),
    @(120,16,'n',2.45),
    @(120,17,'n',0),
    @(120,18,'n',1.975),
    @(120,19,'n',1.825),
    @(120,20,'n',2.5),
    @(120,21,'n',2.025),
    @(120,22,'n',1.775),
    @(120,23,'n',0),
    @(120,24,'n',0),
    @(120,25,'n',0),
    @(120,26,'n',0),
    @(120,27,'n',0),
    @(121,1,'n',119),
    @(121,2,'n',6814420),
    @(121,3,'s','Slovenia Prva Liga'),
    @(121,4,'s','Slovenia Prva Liga'),
    @(121,5,'n',45361.45833333334),
    @(121,6,'s','NS Mura'),
    @(121,7,'s','Olimpija Ljubljana'),
    @(121,11,'n',5.75),
    @(121,12,'n',4),
    @(121,13,'n',1.5),
    @(121,14,'n',5.25),
    @(121,15,'n',3.8),
    @(121,16,'n',1.55),
    @(121,17,'n',1),
    @(121,18,'n',1.775),
    @(121,19,'n',2.025),
    @(121,20,'n',2.5),
    @(121,21,'n',1.85),
    @(121,22,'n',1.95),
    @(121,23,'n',0),
    @(121,24,'n',0),
    @(121,25,'n',0),
    @(121,26,'n',0),
    @(121,27,'n',0),
    @(122,1,'n',120),
    @(122,2,'n',6814422),
    @(122,3,'s','Slovenia Prva Liga'),
    @(122,4,'s','Slovenia Prva Liga'),
    @(122,5,'n',45361.67708333334),
    @(122,6,'s','NK Domzale'),
    @(122,7,'s','FC Koper'),
    @(122,11,'n',3),
    @(122,12,'n',3.25),
    @(122,13,'n',2.2),
    @(122,14,'n',3.1),
    @(122,15,'n',3.25),
    @(122,16,'n',2.15),
    @(122,17,'n',0.25),
    @(122,18,'n',1.9),
    @(122,19,'n',1.9),
    @(122,20,'n',2.5),
    @(122,21,'n',1.95),
    @(122,22,'n',1.85),
    @(122,23,'n',0),
    @(122,24,'n',0),
    @(122,25,'n',0),
    @(122,26,'n',0),
    @(122,27,'n',0)
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $t = $item[2]
    $v = $item[3]
    $cell = $ws.Cells.Item($r, $c)
    if ($t -eq 's') {
        $cell.Value = [string]$v
    } else {
        $cell.Value = [double]$v
    }
}
